$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $NewValue)
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '62.846.21'
Set-TextValue $ws 'E2' '  +0.19%  '
Set-TextValue $ws 'D3' '2.462.87'
Set-TextValue $ws 'E3' '  +0.71%  '
Set-TextValue $ws 'E4' '  +0.05%  '
Set-TextValue $ws 'D5' '574.80'
Set-TextValue $ws 'E5' '  -0.18%  '
Set-TextValue $ws 'D6' '146.96'
Set-TextValue $ws 'E6' '  +0.88%  '
Set-TextValue $ws 'E7' '  -0.02%  '
Set-TextValue $ws 'E8' '  -0.63%  '
Set-TextValue $ws 'D9' '2.462.74'
Set-TextValue $ws 'E9' '  +0.74%  '
Set-TextValue $ws 'E10' '  +0.79%  '
Set-TextValue $ws 'D11' '0.162'
Set-TextValue $ws 'E11' '  -0.47%  '
Set-TextValue $ws 'E12' '  +0.42%  '
Set-TextValue $ws 'D13' '0.357'
Set-TextValue $ws 'E13' '  +1.18%  '
Set-TextValue $ws 'D14' '29.01'
Set-TextValue $ws 'E14' '  +2.95%  '
Set-TextValue $ws 'E15' '  -0.04%  '
Set-TextValue $ws 'D16' '2.909.61'
Set-TextValue $ws 'E16' '  +0.73%  '
Set-TextValue $ws 'D17' '62.755.96'
Set-TextValue $ws 'E17' '  +0.26%  '
Set-TextValue $ws 'D18' '2.460.87'
Set-TextValue $ws 'E18' '  +0.60%  '
Set-TextValue $ws 'D19' '7.94'
Set-TextValue $ws 'E19' '  +0.24%  '
Set-TextValue $ws 'D20' '11.00'
Set-TextValue $ws 'D21' '326.79'
Set-TextValue $ws 'E21' '  -0.82%  '
Set-TextValue $ws 'E22' '  -0.03%  '
Set-TextValue $ws 'E23' '  +8.59%  '
Set-TextValue $ws 'E24' '  -0.05%  '
Set-TextValue $ws 'D25' '10.04'
Set-TextValue $ws 'E25' '  +17.98%  '
Set-TextValue $ws 'D26' '65.46'
Set-TextValue $ws 'E26' '  -1.03%  '
Set-TextValue $ws 'D27' '646.24'
Set-TextValue $ws 'E27' '  -0.01%  '
Set-TextValue $ws 'D28' '0.0₃0988'
Set-TextValue $ws 'E28' '  +0.11%  '
Set-TextValue $ws 'D29' '2.591.67'
Set-TextValue $ws 'D30' '0.999'
Set-TextValue $ws 'E30' '  -15.13%  '
Set-TextValue $ws 'E31' '  -0.60%  '
Set-TextValue $ws 'E32' '  -2.77%  '
Set-TextValue $ws 'E33' '  -1.46%  '
Set-TextValue $ws 'E34' '  -2.88%  '
Set-TextValue $ws 'D35' '0.999'
Set-TextValue $ws 'E35' '  +0.00%  '
Set-TextValue $ws 'E37' '  -0.32%  '
Set-TextValue $ws 'E38' '  +4.79%  '
Set-TextValue $ws 'D39' '0.369'
Set-TextValue $ws 'E39' '  -1.33%  '
Set-TextValue $ws 'B40' 'RenderToken'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws 'D40' '5.40'
Set-TextValue $ws 'E40' '  -1.69%  '
Set-TextValue $ws 'B41' 'Monero'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D41' '151.41'
Set-TextValue $ws 'E41' '  -1.32%  '
Set-TextValue $ws 'B42' 'EthereumClassic'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D42' '18.70'
Set-TextValue $ws 'E42' '  -0.23%  '
Set-TextValue $ws 'E43' '  -0.80%  '
Set-TextValue $ws 'D44' '0.0₆0308'
Set-TextValue $ws 'E44' '  -42.95%  '
Set-TextValue $ws 'E45' '  +0.01%  '
Set-TextValue $ws 'D46' '152.12'
Set-TextValue $ws 'E47' '  +2.11%  '
Set-TextValue $ws 'D48' '3.58'
Set-TextValue $ws 'E48' '  -1.55%  '
Set-TextValue $ws 'D49' '20.55'
Set-TextValue $ws 'E49' '  -0.37%  '
Set-TextValue $ws 'D50' '0.608'
Set-TextValue $ws 'E50' '  +0.44%  '
Set-TextValue $ws 'D51' '0.0511'
Set-TextValue $ws 'E51' '  -1.28%  '
